# Update values in Sheet1 per the target diff.
# (Commit message "Update Name of Algo" refers to the overall batch of result
#  files for different algorithms; for this specific workbook the only
#  change is to the numeric data values below.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Value = -7.808
$ws.Range("C9").Value = -11.445
$ws.Range("C18").Value = -12.314
$ws.Range("C20").Value = -12.581
$ws.Range("E21").Value = 13.146
